# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    existing "2022-Q2" tab) and fill it with the quarter's fund-holdings
#    table.
# 2) Insert a new top row in the "总计" (totals) sheet for the 2022-Q3
#    summary figures, pushing the earlier history down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q3" worksheet
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($anchor)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $q3.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$q3Data = @(
    @(0, "007291", "汇丰晋信港股通双核策略混合", "6.17", "94.56", "9.23", "0.5695", 2),
    @(1, "002332", "汇丰晋信沪港深股票A", "5.45", "94.22", "5.32", "0.2899", 4),
    @(2, "010751", "宝盈优质成长混合A", "4.72", "92.28", "4.02", "0.1897", 9),
    @(3, "011017", "鹏扬景明一年持有期混合", "18.76", "23.56", "0.54", "0.1013", 8),
    @(4, "007107", "太平 MSCI 香港价值增强指数A", "0.97", "90.37", "4.47", "0.0434", 6),
    @(5, "007132", "长城港股通价值精选多策略混合", "0.72", "80.94", "4.84", "0.0348", 5),
    @(6, "010752", "宝盈优质成长混合C", "0.76", "92.28", "4.02", "0.0306", 9),
    @(7, "009130", "鹏扬景恒六个月持有期混合A", "4.75", "24.46", "0.57", "0.0271", 10),
    @(8, "003413", "华泰柏瑞新经济沪港深混合", "0.42", "86.45", "6.09", "0.0256", 3),
    @(9, "002333", "汇丰晋信沪港深股票C", "0.48", "94.22", "5.32", "0.0255", 4),
    @(10, "005255", "浦银安盛港股通量化混合A", "0.29", "78.68", "4.33", "0.0126", 7),
    @(11, "009131", "鹏扬景恒六个月持有期混合C", "1.54", "24.46", "0.57", "0.0088", 10),
    @(12, "011243", "万家惠裕回报6个月持有期混合A", "1.54", "27.67", "0.37", "0.0057", 10),
    @(13, "013224", "浦银安盛港股通量化混合C", "0.05", "78.68", "4.33", "0.0022", 7),
    @(14, "011244", "万家惠裕回报6个月持有期混合C", "0.12", "27.67", "0.37", "0.0004", 10),
    @(15, "007108", "太平 MSCI 香港价值增强指数C", "0.00", "90.37", "4.47", "0", 6)
)

# Columns B-G hold text (fund codes keep leading zeros, figures are stored
# as strings in the source data) - force text formatting before writing so
# values like "007291" / "6.17" are not reinterpreted as numbers.
$q3.Range("B2:G17").NumberFormat = "@"

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# G17 is the one exception: a genuine numeric 0, not text.
$q3.Range("G17").NumberFormat = "General"
$q3.Range("G17").Value = 0

# Match the bold / centered / thin-bordered header-row + index-column look
# used on every other quarter sheet.
$styleRanges = @($q3.Range("B1:H1"), $q3.Range("A2:A17"))
foreach ($sr in $styleRanges) {
    $sr.Font.Bold = $true
    $sr.Borders.LineStyle = 1
    $sr.HorizontalAlignment = -4108
    $sr.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------
# 2. "总计" (totals) sheet: new 2022-Q3 row on top, history shifts down
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows.Item(2).Insert()
$totals.Rows.Item(2).ClearFormats()

# Re-use the existing index-column style (now on row 3) for the new row.
$totals.Range("A3").Copy() | Out-Null
$totals.Range("A2").PasteSpecial(-4122) | Out-Null

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 16
$totals.Range("D2").Value = 1.37
